$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change status ("Im Besitz") of episodes (rows 2-6) from "Nein" to "Ja"
$ws.Range("C2:C6").Value = "Ja"

# Update the selection to match row 7 being selected (active cell A7)
$ws.Range("A7:XFD7").Select()
